# The "Peak Annotation Files" sheet previously stored the full relative
# repo path to each peak annotation file in column A.  The fix changes
# these to just the bare file name (the loader now resolves the
# containing directory itself), matching the file names already used
# elsewhere in the workbook (e.g. "Sequences"."Peak Annotation File Name").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peak Annotation Files")

$ws.Range("A2").Value = "accucor1.xlsx"
$ws.Range("A3").Value = "accucor2.xlsx"
